$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.943.17'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.203.71'
$ws.Range("E3").Value = '  -0.92%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '293.95'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = '86.90'
$ws.Range("E6").Value = '  -1.58%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").Value = '51.41'
$ws.Range("E10").Value = '  +7.15%  '
$ws.Range("D11").Value = '30.50'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '0.0778'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E13").Value = '  +2.70%  '
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '2.546.18'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '13.75'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '2.206.71'
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '39.860.62'
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("D20").Value = '0.0₃0885'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").Value = '11.16'
$ws.Range("E21").Value = '  -3.54%  '
$ws.Range("E22").Value = '  -1.48%  '
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").Value = '234.94'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").Value = '2.47'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D28").Value = '23.01'
$ws.Range("E28").Value = '  +1.44%  '
$ws.Range("D29").Value = '2.08'
$ws.Range("E29").Value = '  -1.22%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").Value = '159.32'
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("D32").Value = '31.63'
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '4.94'
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '3.02'
$ws.Range("E35").Value = '  +5.21%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("D39").Value = '0.0997'
$ws.Range("E39").Value = '  +0.97%  '
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("D41").Value = '15.44'
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("D42").Value = '2.069.45'
$ws.Range("E42").Value = '  -2.15%  '
$ws.Range("D43").Value = '3.72'
$ws.Range("E43").Value = '  -3.23%  '
$ws.Range("D44").Value = '19.06'
$ws.Range("E44").Value = '  +6.88%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '9.84'
$ws.Range("E46").Value = '  -1.25%  '
$ws.Range("E47").Value = '  +2.41%  '
$ws.Range("E48").Value = '  -9.32%  '
$ws.Range("D49").Value = '2.420.11'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  +1.47%  '
$ws.Range("E51").Value = '  -0.31%  '
